# Daily attendance processing - re-sort the "Recorded By" (column G) list of
# names/emails alphabetically (ordinal, case-sensitive, ascending) for every
# data row of the active worksheet.

function Compare-Ordinal($s1, $s2) {
    $len1 = $s1.Length
    $len2 = $s2.Length
    $minLen = [Math]::Min($len1, $len2)
    for ($ci = 0; $ci -lt $minLen; $ci++) {
        $c1 = [int][char]$s1[$ci]
        $c2 = [int][char]$s2[$ci]
        if ($c1 -lt $c2) { return -1 }
        if ($c1 -gt $c2) { return 1 }
    }
    if ($len1 -lt $len2) { return -1 }
    if ($len1 -gt $len2) { return 1 }
    return 0
}

function Sort-Ordinal($items) {
    $arr = @($items)
    $n = $arr.Count
    for ($si = 0; $si -lt $n; $si++) {
        for ($sj = 0; $sj -lt ($n - $si - 1); $sj++) {
            $cmp = Compare-Ordinal $arr[$sj] $arr[$sj + 1]
            if ($cmp -gt 0) {
                $tmp = $arr[$sj]
                $arr[$sj] = $arr[$sj + 1]
                $arr[$sj + 1] = $tmp
            }
        }
    }
    return $arr
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $val = $cell.Text
    if ($val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -gt 1) {
            $sorted = Sort-Ordinal $parts
            $joined = $sorted -join ", "
            if ($joined -ne $val) {
                $cell.Value = $joined
            }
        }
    }
}
